# Auto-generated edit script: updates numeric columns H-N for specific Leve rows
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5778.1577  # H86: was 6096.3887
$ws.Cells.Item(86, 10).Value = 14618.143  # J86: was 17046.166
$ws.Cells.Item(86, 12).Value = 14618.143  # L86: was 17046.166
$ws.Cells.Item(86, 14).Value = -16864.143  # N86: was -19292.166

$ws.Cells.Item(89, 8).Value = 5778.1577  # H89: was 6096.3887
$ws.Cells.Item(89, 10).Value = 14618.143  # J89: was 17046.166
$ws.Cells.Item(89, 12).Value = 73090.715  # L89: was 85230.83
$ws.Cells.Item(89, 14).Value = -84322.715  # N89: was -96462.83

$ws.Cells.Item(113, 8).Value = 35718180  # H113: was 34486600
$ws.Cells.Item(113, 9).Value = 45457790  # I113: was 47622372
$ws.Cells.Item(113, 10).Value = 6268.6665  # J113: was 5201.5
$ws.Cells.Item(113, 11).Value = 45457790  # K113: was 47622372
$ws.Cells.Item(113, 12).Value = 6268.6665  # L113: was 5201.5
$ws.Cells.Item(113, 13).Value = -45454536  # M113: was -47619118
$ws.Cells.Item(113, 14).Value = -12776.6665  # N113: was -11709.5

$ws.Cells.Item(129, 8).Value = 251171.03  # H129: was 313744.03
$ws.Cells.Item(129, 10).Value = 287005.16  # J129: was 371783.28
$ws.Cells.Item(129, 12).Value = 861015.48  # L129: was 1115349.84
$ws.Cells.Item(129, 14).Value = -871015.48  # N129: was -1125349.84

$ws.Cells.Item(137, 8).Value = 4588  # H137: was 4489.3335
$ws.Cells.Item(137, 9).Value = 4950.1665  # I137: was 5366.8335
$ws.Cells.Item(137, 10).Value = 3501.5  # J137: was 2734.3333
$ws.Cells.Item(137, 11).Value = 14850.4995  # K137: was 16100.5005
$ws.Cells.Item(137, 12).Value = 10504.5  # L137: was 8202.999899999999
$ws.Cells.Item(137, 13).Value = -12300.4995  # M137: was -13550.5005
$ws.Cells.Item(137, 14).Value = -15604.5  # N137: was -13302.9999

$ws.Cells.Item(138, 8).Value = 27781246  # H138: was 26319252
$ws.Cells.Item(138, 10).Value = 4376.905  # J138: was 4290.2173
$ws.Cells.Item(138, 12).Value = 13130.715  # L138: was 12870.6519
$ws.Cells.Item(138, 14).Value = -23410.715  # N138: was -23150.6519

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 361247.12  # H61: was 368635.5
$ws.Cells.Item(61, 9).Value = 392270.8  # I61: was 400916.47
$ws.Cells.Item(61, 10).Value = 4474.75  # J61: was 5474.75
$ws.Cells.Item(61, 11).Value = 392270.8  # K61: was 400916.47
$ws.Cells.Item(61, 12).Value = 4474.75  # L61: was 5474.75
$ws.Cells.Item(61, 13).Value = -392058.8  # M61: was -400704.47
$ws.Cells.Item(61, 14).Value = -4898.75  # N61: was -5898.75

$ws.Cells.Item(74, 8).Value = 55558576  # H74: was 37039150
$ws.Cells.Item(74, 9).Value = 58826548  # I74: was 38463620
$ws.Cells.Item(74, 11).Value = 58826548  # K74: was 38463620
$ws.Cells.Item(74, 13).Value = -58825674  # M74: was -38462746

$ws.Cells.Item(77, 8).Value = 55558576  # H77: was 37039150
$ws.Cells.Item(77, 9).Value = 58826548  # I77: was 38463620
$ws.Cells.Item(77, 11).Value = 294132740  # K77: was 192318100
$ws.Cells.Item(77, 13).Value = -294128372  # M77: was -192313732

$ws.Cells.Item(97, 8).Value = 1878.8182  # H97: was 1926.7
$ws.Cells.Item(97, 9).Value = 1729.6666  # I97: was 1813.7646
$ws.Cells.Item(97, 10).Value = 2550  # J97: was 2566.6667
$ws.Cells.Item(97, 11).Value = 1729.6666  # K97: was 1813.7646
$ws.Cells.Item(97, 12).Value = 2550  # L97: was 2566.6667
$ws.Cells.Item(97, 13).Value = -1233.6666  # M97: was -1317.7646
$ws.Cells.Item(97, 14).Value = -3542  # N97: was -3558.6667

$ws.Cells.Item(102, 8).Value = 1001.82355  # H102: was 1057.4615
$ws.Cells.Item(102, 9).Value = 939.5  # I102: was 1057.4615
$ws.Cells.Item(102, 10).Value = 1999  # J102: was 0
$ws.Cells.Item(102, 11).Value = 939.5  # K102: was 1057.4615
$ws.Cells.Item(102, 12).Value = 1999  # L102: was 0
$ws.Cells.Item(102, 13).Value = 682.5  # M102: was 564.5385000000001
$ws.Cells.Item(102, 14).Value = -5243  # N102: was None

$ws.Cells.Item(132, 8).Value = 17188.273  # H132: was 18255.033
$ws.Cells.Item(132, 9).Value = 2221.4546  # I132: was 2378
$ws.Cells.Item(132, 10).Value = 47121.91  # J132: was 47122.363
$ws.Cells.Item(132, 11).Value = 6664.3638  # K132: was 7134
$ws.Cells.Item(132, 12).Value = 141365.73  # L132: was 141367.089
$ws.Cells.Item(132, 13).Value = -4134.3638  # M132: was -4604
$ws.Cells.Item(132, 14).Value = -146425.73  # N132: was -146427.089

$ws.Cells.Item(136, 8).Value = 361247.12  # H136: was 368635.5
$ws.Cells.Item(136, 9).Value = 392270.8  # I136: was 400916.47
$ws.Cells.Item(136, 10).Value = 4474.75  # J136: was 5474.75
$ws.Cells.Item(136, 11).Value = 1176812.4  # K136: was 1202749.41
$ws.Cells.Item(136, 12).Value = 13424.25  # L136: was 16424.25
$ws.Cells.Item(136, 13).Value = -1174262.4  # M136: was -1200199.41
$ws.Cells.Item(136, 14).Value = -18524.25  # N136: was -21524.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 548.0476  # H94: was 571.5
$ws.Cells.Item(94, 9).Value = 473.35294  # I94: was 498
$ws.Cells.Item(94, 11).Value = 473.35294  # K94: was 498
$ws.Cells.Item(94, 13).Value = -22.35293999999999  # M94: was -47

$ws.Cells.Item(134, 8).Value = 2678.3696  # H134: was 2772.2046
$ws.Cells.Item(134, 9).Value = 2886.375  # I134: was 2886.925
$ws.Cells.Item(134, 10).Value = 1291.6666  # J134: was 1625
$ws.Cells.Item(134, 11).Value = 8659.125  # K134: was 8660.775000000001
$ws.Cells.Item(134, 12).Value = 3874.9998  # L134: was 4875
$ws.Cells.Item(134, 13).Value = -6124.125  # M134: was -6125.775000000001
$ws.Cells.Item(134, 14).Value = -8944.9998  # N134: was -9945

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5799.769  # H31: was 4494.4
$ws.Cells.Item(31, 9).Value = 4067.7856  # I31: was 2852.682
$ws.Cells.Item(31, 10).Value = 7820.4165  # J31: was 7272.6924
$ws.Cells.Item(31, 11).Value = 4067.7856  # K31: was 2852.682
$ws.Cells.Item(31, 12).Value = 7820.4165  # L31: was 7272.6924
$ws.Cells.Item(31, 13).Value = -3772.7856  # M31: was -2557.682
$ws.Cells.Item(31, 14).Value = -8410.416499999999  # N31: was -7862.6924

$ws.Cells.Item(34, 8).Value = 5799.769  # H34: was 4494.4
$ws.Cells.Item(34, 9).Value = 4067.7856  # I34: was 2852.682
$ws.Cells.Item(34, 10).Value = 7820.4165  # J34: was 7272.6924
$ws.Cells.Item(34, 11).Value = 4067.7856  # K34: was 2852.682
$ws.Cells.Item(34, 12).Value = 7820.4165  # L34: was 7272.6924
$ws.Cells.Item(34, 13).Value = -3865.7856  # M34: was -2650.682
$ws.Cells.Item(34, 14).Value = -8224.416499999999  # N34: was -7676.6924

$ws.Cells.Item(58, 8).Value = 11246.225  # H58: was 12768.628
$ws.Cells.Item(58, 9).Value = 903.7353000000001  # I58: was 965.70966
$ws.Cells.Item(58, 10).Value = 34689.2  # J58: was 43259.5
$ws.Cells.Item(58, 11).Value = 903.7353000000001  # K58: was 965.70966
$ws.Cells.Item(58, 12).Value = 34689.2  # L58: was 43259.5
$ws.Cells.Item(58, 13).Value = -700.7353000000001  # M58: was -762.70966
$ws.Cells.Item(58, 14).Value = -35095.2  # N58: was -43665.5

$ws.Cells.Item(136, 8).Value = 11246.225  # H136: was 12768.628
$ws.Cells.Item(136, 9).Value = 903.7353000000001  # I136: was 965.70966
$ws.Cells.Item(136, 10).Value = 34689.2  # J136: was 43259.5
$ws.Cells.Item(136, 11).Value = 2711.2059  # K136: was 2897.12898
$ws.Cells.Item(136, 12).Value = 104067.6  # L136: was 129778.5
$ws.Cells.Item(136, 13).Value = -161.2058999999999  # M136: was -347.12898
$ws.Cells.Item(136, 14).Value = -109167.6  # N136: was -134878.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1122.4117  # H5: was 1236.7333
$ws.Cells.Item(5, 9).Value = 462.2  # I5: was 511.5
$ws.Cells.Item(5, 11).Value = 1386.6  # K5: was 1534.5
$ws.Cells.Item(5, 13).Value = -1274.6  # M5: was -1422.5

$ws.Cells.Item(118, 8).Value = 55558336  # H118: was 55558348
$ws.Cells.Item(118, 10).Value = 4819.6  # J118: was 4839.6
$ws.Cells.Item(118, 12).Value = 14458.8  # L118: was 14518.8
$ws.Cells.Item(118, 14).Value = -16944.8  # N118: was -17004.8

$ws.Cells.Item(131, 8).Value = 776.01  # H131: was 102841.13
$ws.Cells.Item(131, 9).Value = 466.66666  # I131: was 607.5
$ws.Cells.Item(131, 10).Value = 795.7553  # J131: was 111928.57
$ws.Cells.Item(131, 11).Value = 1399.99998  # K131: was 1822.5
$ws.Cells.Item(131, 12).Value = 2387.2659  # L131: was 335785.71
$ws.Cells.Item(131, 13).Value = 3640.00002  # M131: was 3217.5
$ws.Cells.Item(131, 14).Value = -12467.2659  # N131: was -345865.71

$ws.Cells.Item(135, 8).Value = 1122.4117  # H135: was 1236.7333
$ws.Cells.Item(135, 9).Value = 462.2  # I135: was 511.5
$ws.Cells.Item(135, 11).Value = 4159.8  # K135: was 4603.5
$ws.Cells.Item(135, 13).Value = -1624.8  # M135: was -2068.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 17950  # H70: was 18946.154
$ws.Cells.Item(70, 9).Value = 19483.334  # I70: was 22380
$ws.Cells.Item(70, 11).Value = 19483.334  # K70: was 22380
$ws.Cells.Item(70, 13).Value = -19213.334  # M70: was -22110

$ws.Cells.Item(73, 8).Value = 17950  # H73: was 18946.154
$ws.Cells.Item(73, 9).Value = 19483.334  # I73: was 22380
$ws.Cells.Item(73, 11).Value = 19483.334  # K73: was 22380
$ws.Cells.Item(73, 13).Value = -18547.334  # M73: was -21444

$ws.Cells.Item(97, 8).Value = 1057.6086  # H97: was 1205.2273
$ws.Cells.Item(97, 9).Value = 1107.45  # I97: was 1226.0555
$ws.Cells.Item(97, 10).Value = 725.3333  # J97: was 1111.5
$ws.Cells.Item(97, 11).Value = 1107.45  # K97: was 1226.0555
$ws.Cells.Item(97, 12).Value = 725.3333  # L97: was 1111.5
$ws.Cells.Item(97, 13).Value = -611.45  # M97: was -730.0554999999999
$ws.Cells.Item(97, 14).Value = -1717.3333  # N97: was -2103.5

$ws.Cells.Item(102, 8).Value = 20835274  # H102: was 22729336
$ws.Cells.Item(102, 9).Value = 22728866  # I102: was 25001692
$ws.Cells.Item(102, 11).Value = 22728866  # K102: was 25001692
$ws.Cells.Item(102, 13).Value = -22727244  # M102: was -25000070

$ws.Cells.Item(107, 8).Value = 5494701  # H107: was 6410472.5
$ws.Cells.Item(107, 9).Value = 231.90909  # I107: was 250.1
$ws.Cells.Item(107, 10).Value = 25641086  # J107: was 38461584
$ws.Cells.Item(107, 11).Value = 231.90909  # K107: was 250.1
$ws.Cells.Item(107, 12).Value = 25641086  # L107: was 38461584
$ws.Cells.Item(107, 13).Value = 1688.09091  # M107: was 1669.9
$ws.Cells.Item(107, 14).Value = -25644926  # N107: was -38465424

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1781.8334  # H68: was 1850.5652
$ws.Cells.Item(68, 9).Value = 1839.4286  # I68: was 1965.4615
$ws.Cells.Item(68, 11).Value = 1839.4286  # K68: was 1965.4615
$ws.Cells.Item(68, 13).Value = -1090.4286  # M68: was -1216.4615

$ws.Cells.Item(71, 8).Value = 1781.8334  # H71: was 1850.5652
$ws.Cells.Item(71, 9).Value = 1839.4286  # I71: was 1965.4615
$ws.Cells.Item(71, 11).Value = 9197.143  # K71: was 9827.307499999999
$ws.Cells.Item(71, 13).Value = -5453.143  # M71: was -6083.307499999999

$ws.Cells.Item(132, 8).Value = 1171.2924  # H132: was 1306.9454
$ws.Cells.Item(132, 9).Value = 1059.0944  # I132: was 1159.6522
$ws.Cells.Item(132, 10).Value = 1666.8334  # J132: was 2059.7778
$ws.Cells.Item(132, 11).Value = 3177.2832  # K132: was 3478.9566
$ws.Cells.Item(132, 12).Value = 5000.5002  # L132: was 6179.3334
$ws.Cells.Item(132, 13).Value = -647.2831999999999  # M132: was -948.9566
$ws.Cells.Item(132, 14).Value = -10060.5002  # N132: was -11239.3334

$ws.Cells.Item(136, 8).Value = 1144.7142  # H136: was 742.1081
$ws.Cells.Item(136, 9).Value = 1144.7142  # I136: was 759.94446
$ws.Cells.Item(136, 10).Value = 0  # J136: was 100
$ws.Cells.Item(136, 11).Value = 3434.1426  # K136: was 2279.83338
$ws.Cells.Item(136, 12).Value = 0  # L136: was 300
$ws.Cells.Item(136, 13).Value = -884.1425999999997  # M136: was 270.16662
$ws.Cells.Item(136, 14).ClearContents()  # N136: was -5400

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 757.6923  # H132: was 835.37933
$ws.Cells.Item(132, 9).Value = 534.2075  # I132: was 578.4167
$ws.Cells.Item(132, 10).Value = 1744.75  # J132: was 2068.8
$ws.Cells.Item(132, 11).Value = 1602.6225  # K132: was 1735.2501
$ws.Cells.Item(132, 12).Value = 5234.25  # L132: was 6206.400000000001
$ws.Cells.Item(132, 13).Value = 927.3775000000001  # M132: was 794.7499
$ws.Cells.Item(132, 14).Value = -10294.25  # N132: was -11266.4

$ws.Cells.Item(136, 8).Value = 16668705  # H136: was 16131008
$ws.Cells.Item(136, 9).Value = 27778844  # I136: was 25642020
$ws.Cells.Item(136, 10).Value = 3495.625  # J136: was 3638.6956
$ws.Cells.Item(136, 11).Value = 83336532  # K136: was 76926060
$ws.Cells.Item(136, 12).Value = 10486.875  # L136: was 10916.0868
$ws.Cells.Item(136, 13).Value = -83333982  # M136: was -76923510
$ws.Cells.Item(136, 14).Value = -15586.875  # N136: was -16016.0868
